$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F2 (Measles query): replace "/" separators with " OR ", drop stray ",\n-fake" tail ---
$ws.Range("F2").Value = "measles OR sarampion OR rougeole OR sarampo OR gafeira OR morrinha"

# --- F3 (COVID-19 query): replace "/" separators with " OR " ---
$ws.Range("F3").Value = 'coronavirus OR "novel coronavirus" OR ncov OR "2019-ncov" OR covid-19 OR sars-covid-2 OR "nuevo coronavirus" OR "nouveau coronavirus" OR "novo coronavirus"'

# --- F4 (COVID-19 outbreaks query): rewrite using AND / OR, keep "/" for synonym lists ---
$cell = $ws.Range("F4")
$newText = "coronavirus AND outbreak/cluster/school OR covid-19 AND outbreak/cluster/school`n-fake"
$cell.Value = $newText

# Recreate the five rich-text runs at their original boundaries:
#   1: "coronavirus AND "            (1-16)
#   2: "outbreak/cluster/school"     (17-39)
#   3: " OR covid-19 AND "           (40-56)
#   4: "outbreak/cluster/school\n"   (57-80)
#   5: "-fake"                       (81-85)
# Runs 1/3/5 get an explicit (visually no-op) color so the engine keeps them
# as distinct runs instead of collapsing back into one.
$cell.Characters(1, 16).Font.Color = 0
$cell.Characters(40, 17).Font.Color = 0
$cell.Characters(81, 5).Font.Color = 0

# --- Selection / autosave housekeeping ---
$ws.Range("F7").Select()
